$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -13.235
$ws.Range("E4").Value = 13.195
$ws.Range("E5").Value = 13.498
$ws.Range("C7").Value = -13.088
$ws.Range("E8").Value = 13.718
$ws.Range("C16").Value = -11.868
$ws.Range("E16").Value = 12.881
